$d = $word.ActiveDocument

function Safe-Replace($doc, $oldText, $newText) {
    $rng = $doc.Content
    $found = $rng.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Output "NOT FOUND: $oldText"
        return
    }
    $targetStart = $rng.Start
    $targetEnd = $rng.End

    $full = $doc.Content
    $docStart = $full.Start
    $docEnd = $full.End

    $didBefore = $false
    $didAfter = $false

    # Shield only a minimal sliver immediately before/after the target so that
    # the run-coalescing pass triggered by the text edit cannot fuse the
    # (untouched) neighboring runs into the edited run OR into each other.
    if ($targetStart -gt $docStart) {
        $before = $doc.Range($targetStart - 1, $targetStart)
        $before.Font.Bold = 1
        $didBefore = $true
    }
    if ($targetEnd -lt $docEnd) {
        $after = $doc.Range($targetEnd, $targetEnd + 1)
        $after.Font.Bold = 1
        $didAfter = $true
    }

    $target = $doc.Range($targetStart, $targetEnd)
    $target.Text = $newText

    $newTargetEnd = $targetStart + $newText.Length

    if ($didBefore) {
        $before2 = $doc.Range($targetStart - 1, $targetStart)
        $before2.Font.Bold = 0
    }
    if ($didAfter) {
        $after2 = $doc.Range($newTargetEnd, $newTargetEnd + 1)
        $after2.Font.Bold = 0
    }
}

Safe-Replace $d "Unveiling the Complexity of Cosmic Phenomena" "Exploring the Marvelous World of Chemistry: A Journey into the Realm of Elements and Compounds"
Safe-Replace $d " Neil deGrasse Tyson" " Eleanor Stanton"
Safe-Replace $d "NeilTyson@SpaceInstitute" "estanton@edu"

Safe-Replace $d "Embarking on an enchanting voyage to unravel the complexities of cosmic phenomena, we delve into the profound mysteries that enchant our universe" "Chemistry, an intriguing and impactful science, unveils the hidden intricacies of matter and its diverse interactions"

Safe-Replace $d " From the grand tapestry of galaxies that stretch across unfathomable distances, to the enigmatic fabric of space and time warping around celestial bodies, our quest for comprehension leads us down a path of captivating discoveries" " From the vast universe to the microscopic realm within our bodies, chemistry plays a pivotal role in shaping our world"

Safe-Replace $d " Through meticulous observation, tireless calculations, and imaginative leap, humanity continues to unlock the secrets of the cosmos, redefining our perception of existence itself" " As we delve into the fascinating tapestry of chemistry, we embark on an exhilarating exploration of the elements that constitute everything around us and the myriad compounds formed through their intricate combinations"

Safe-Replace $d "As our telescopes peer deeper into the vast expanse, we encounter distant worlds that ignite our curiosity and contemplation" "In this realm of substances, we uncover the fundamental principles governing chemical reactions, witnessing the wondrous transformations of matter into new entities with unique properties"

Safe-Replace $d " The interplay of cosmic forces, the birth and death of stars, and the symphony of interactions between celestial bodies fuel our inquiry into the fundamental principles that govern the universe's evolution" " The symphony of chemistry encompasses myriad concepts, from atomic structures and bonding arrangements to energy transfer and reaction dynamics"

Safe-Replace $d " Each cosmic event, each intricate celestial dance, holds clues to unraveling mysteries that have captivated humankind for millennia, painting a breathtaking tapestry of cosmic beauty and awe" " Each element, with its distinctive characteristics, contributes to the intricate dance of chemical interactions, orchestrating the formation of countless compounds with diverse applications in fields ranging from medicine to materials science"

Safe-Replace $d "Yet, the complexities of the universe extend beyond the reaches of our tangible world, delving into realms that transcend our current understanding" "As we unravel the enigmas of chemistry, we gain invaluable insights into the natural world, unveiling the intricate mechanisms underlying life itself"

Safe-Replace $d " From the mysteries of dark matter and energy, whose enigmatic nature eludes our grasp, to the theoretical concept of multiple universes and the complexities of multi-dimensional space, our exploration into the cosmos opens doors to realms that challenge our conceptual boundaries" " From the intricate workings of photosynthesis, the process by which plants convert sunlight into energy, to the intricate pathways of cellular respiration, the fundamental energy-generating process within living organisms, chemistry unveils the symphony of life at its most fundamental level"

Safe-Replace $d " In this journey of seeking cosmic knowledge, we embark on an intellectual adventure where wonder and enigma intersect, inviting us to the depths of the universe's profound secrets" " Its principles permeate every aspect of our existence, shaping the materials we use, the medicines that heal us, and the intricate complexity of the living world"

Safe-Replace $d "Our journey of exploration into cosmic phenomena unveils a harmonious tapestry of elegance, mystery, and boundless beauty" "This essay embarks on an enthralling exploration of chemistry, venturing into the captivating realm of elements, compounds, and their captivating interactions"

Safe-Replace $d " As we continue to probe the cosmos with unrelenting curiosity, we are reminded that the path to cosmic understanding is a never-ending adventure, filled with captivating discoveries and awe-inspiring revelations that paint a mesmerizing portrait of the intricate universe we inhabit" " Through the study of chemistry, we gain a profound understanding of the natural world, unlocking the secrets of matter and its remarkable transformations, revealing the symphony of life at its most fundamental level"

Write-Output "done basic replacements"
